$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value = 5056287
$ws.Cells.Item(28, 9).Value = 8555731
$ws.Cells.Item(28, 10).Value = 1534.3334
$ws.Cells.Item(28, 11).Value = 8555731
$ws.Cells.Item(28, 12).Value = 1534.3334
$ws.Cells.Item(28, 13).Value = -8555246
$ws.Cells.Item(28, 14).Value = -2504.3334
# Row 33
$ws.Cells.Item(33, 8).Value = 6550
$ws.Cells.Item(33, 9).Value = 100
$ws.Cells.Item(33, 10).Value = 6953.125
$ws.Cells.Item(33, 11).Value = 100
$ws.Cells.Item(33, 12).Value = 6953.125
$ws.Cells.Item(33, 13).Value = 129
$ws.Cells.Item(33, 14).Value = -7411.125
# Row 64
$ws.Cells.Item(64, 8).Value = 3291.9412
$ws.Cells.Item(64, 9).Value = 3180.9092
$ws.Cells.Item(64, 10).Value = 3495.5
$ws.Cells.Item(64, 11).Value = 3180.9092
$ws.Cells.Item(64, 12).Value = 3495.5
$ws.Cells.Item(64, 13).Value = -2932.9092
$ws.Cells.Item(64, 14).Value = -3991.5
# Row 67
$ws.Cells.Item(67, 8).Value = 3291.9412
$ws.Cells.Item(67, 9).Value = 3180.9092
$ws.Cells.Item(67, 10).Value = 3495.5
$ws.Cells.Item(67, 11).Value = 3180.9092
$ws.Cells.Item(67, 12).Value = 3495.5
$ws.Cells.Item(67, 13).Value = -2322.9092
$ws.Cells.Item(67, 14).Value = -5211.5

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 1764.7142
$ws.Cells.Item(61, 9).Value = 2370.3333
$ws.Cells.Item(61, 10).Value = 1310.5
$ws.Cells.Item(61, 11).Value = 2370.3333
$ws.Cells.Item(61, 12).Value = 1310.5
$ws.Cells.Item(61, 13).Value = -2158.3333
$ws.Cells.Item(61, 14).Value = -1734.5
# Row 64
$ws.Cells.Item(64, 8).Value = 100000
$ws.Cells.Item(64, 9).Value = 100000
$ws.Cells.Item(64, 11).Value = 100000
$ws.Cells.Item(64, 13).Value = -99752
# Row 67
$ws.Cells.Item(67, 8).Value = 100000
$ws.Cells.Item(67, 9).Value = 100000
$ws.Cells.Item(67, 11).Value = 100000
$ws.Cells.Item(67, 13).Value = -99142
# Row 132
$ws.Cells.Item(132, 8).Value = 23280904
$ws.Cells.Item(132, 9).Value = 32259620
$ws.Cells.Item(132, 10).Value = 85885.664
$ws.Cells.Item(132, 11).Value = 96778860
$ws.Cells.Item(132, 12).Value = 257656.992
$ws.Cells.Item(132, 13).Value = -96776330
$ws.Cells.Item(132, 14).Value = -262716.992
# Row 136
$ws.Cells.Item(136, 8).Value = 1764.7142
$ws.Cells.Item(136, 9).Value = 2370.3333
$ws.Cells.Item(136, 10).Value = 1310.5
$ws.Cells.Item(136, 11).Value = 7110.999899999999
$ws.Cells.Item(136, 12).Value = 3931.5
$ws.Cells.Item(136, 13).Value = -4560.999899999999
$ws.Cells.Item(136, 14).Value = -9031.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 7409.2856
$ws.Cells.Item(20, 9).Value = 6912.1875
$ws.Cells.Item(20, 10).Value = 9000
$ws.Cells.Item(20, 11).Value = 6912.1875
$ws.Cells.Item(20, 12).Value = 9000
$ws.Cells.Item(20, 13).Value = -6665.1875
$ws.Cells.Item(20, 14).Value = -9494
# Row 132
$ws.Cells.Item(132, 8).Value = 20000
$ws.Cells.Item(132, 10).Value = 20000
$ws.Cells.Item(132, 12).Value = 20000
$ws.Cells.Item(132, 14).Value = -30120

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 14712689
$ws.Cells.Item(31, 9).Value = 17859536
$ws.Cells.Item(31, 10).Value = 27400
$ws.Cells.Item(31, 11).Value = 17859536
$ws.Cells.Item(31, 12).Value = 27400
$ws.Cells.Item(31, 13).Value = -17859241
$ws.Cells.Item(31, 14).Value = -27990
# Row 34
$ws.Cells.Item(34, 8).Value = 14712689
$ws.Cells.Item(34, 9).Value = 17859536
$ws.Cells.Item(34, 10).Value = 27400
$ws.Cells.Item(34, 11).Value = 17859536
$ws.Cells.Item(34, 12).Value = 27400
$ws.Cells.Item(34, 13).Value = -17859334
$ws.Cells.Item(34, 14).Value = -27804
# Row 58
$ws.Cells.Item(58, 8).Value = 1301.5
$ws.Cells.Item(58, 9).Value = 1301.5
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 1301.5
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = -1098.5
$ws.Cells.Item(58, 14).Value = $null
# Row 99
$ws.Cells.Item(99, 8).Value = 247820.12
$ws.Cells.Item(99, 9).Value = 370930.8
$ws.Cells.Item(99, 10).Value = 1598.7391
$ws.Cells.Item(99, 11).Value = 370930.8
$ws.Cells.Item(99, 12).Value = 1598.7391
$ws.Cells.Item(99, 13).Value = -369432.8
$ws.Cells.Item(99, 14).Value = -4594.7391
# Row 126
$ws.Cells.Item(126, 8).Value = 247820.12
$ws.Cells.Item(126, 9).Value = 370930.8
$ws.Cells.Item(126, 10).Value = 1598.7391
$ws.Cells.Item(126, 11).Value = 1112792.4
$ws.Cells.Item(126, 12).Value = 4796.2173
$ws.Cells.Item(126, 13).Value = -1110322.4
$ws.Cells.Item(126, 14).Value = -9736.2173
# Row 132
$ws.Cells.Item(132, 8).Value = 55742
$ws.Cells.Item(132, 9).Value = 2084.6155
$ws.Cells.Item(132, 11).Value = 6253.8465
$ws.Cells.Item(132, 13).Value = -3723.8465
# Row 134
$ws.Cells.Item(134, 8).Value = 1842.2759
$ws.Cells.Item(134, 9).Value = 1376.9474
$ws.Cells.Item(134, 10).Value = 2726.4
$ws.Cells.Item(134, 11).Value = 4130.8422
$ws.Cells.Item(134, 12).Value = 8179.200000000001
$ws.Cells.Item(134, 13).Value = -1595.8422
$ws.Cells.Item(134, 14).Value = -13249.2
# Row 136
$ws.Cells.Item(136, 8).Value = 1301.5
$ws.Cells.Item(136, 9).Value = 1301.5
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 3904.5
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -1354.5
$ws.Cells.Item(136, 14).Value = $null

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Cells.Item(14, 8).Value = 179.18182
$ws.Cells.Item(14, 9).Value = 179.18182
$ws.Cells.Item(14, 11).Value = 537.5454599999999
$ws.Cells.Item(14, 13).Value = -364.5454599999999
# Row 125
$ws.Cells.Item(125, 8).Value = 2475
$ws.Cells.Item(125, 9).Value = 3725
$ws.Cells.Item(125, 10).Value = 2058.3333
$ws.Cells.Item(125, 11).Value = 11175
$ws.Cells.Item(125, 12).Value = 6174.999899999999
$ws.Cells.Item(125, 13).Value = -6255
$ws.Cells.Item(125, 14).Value = -16014.9999
# Row 132
$ws.Cells.Item(132, 8).Value = 1906.6666
$ws.Cells.Item(132, 9).Value = 761.4
$ws.Cells.Item(132, 10).Value = 2479.3
$ws.Cells.Item(132, 11).Value = 6852.599999999999
$ws.Cells.Item(132, 12).Value = 22313.7
$ws.Cells.Item(132, 13).Value = -4322.599999999999
$ws.Cells.Item(132, 14).Value = -27373.7

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 31259160
$ws.Cells.Item(70, 9).Value = 62507988
$ws.Cells.Item(70, 10).Value = 10333.333
$ws.Cells.Item(70, 11).Value = 62507988
$ws.Cells.Item(70, 12).Value = 10333.333
$ws.Cells.Item(70, 13).Value = -62507718
$ws.Cells.Item(70, 14).Value = -10873.333
# Row 73
$ws.Cells.Item(73, 8).Value = 31259160
$ws.Cells.Item(73, 9).Value = 62507988
$ws.Cells.Item(73, 10).Value = 10333.333
$ws.Cells.Item(73, 11).Value = 62507988
$ws.Cells.Item(73, 12).Value = 10333.333
$ws.Cells.Item(73, 13).Value = -62507052
$ws.Cells.Item(73, 14).Value = -12205.333
# Row 96
$ws.Cells.Item(96, 8).Value = 23000
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 10).Value = 23000
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 12).Value = 23000
$ws.Cells.Item(96, 13).Value = $null
$ws.Cells.Item(96, 14).Value = -28492
# Row 132
$ws.Cells.Item(132, 8).Value = 394692.56
$ws.Cells.Item(132, 9).Value = 59944.234
$ws.Cells.Item(132, 10).Value = 912030.9399999999
$ws.Cells.Item(132, 11).Value = 179832.702
$ws.Cells.Item(132, 12).Value = 2736092.82
$ws.Cells.Item(132, 13).Value = -177302.702
$ws.Cells.Item(132, 14).Value = -2741152.82

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Cells.Item(93, 8).Value = 3716.5833
$ws.Cells.Item(93, 9).Value = 4089.5557
$ws.Cells.Item(93, 10).Value = 2597.6667
$ws.Cells.Item(93, 11).Value = 4089.5557
$ws.Cells.Item(93, 12).Value = 2597.6667
$ws.Cells.Item(93, 13).Value = -2841.5557
$ws.Cells.Item(93, 14).Value = -5093.6667
# Row 122
$ws.Cells.Item(122, 8).Value = 2238.7307
$ws.Cells.Item(122, 9).Value = 2011.5555
$ws.Cells.Item(122, 10).Value = 2749.875
$ws.Cells.Item(122, 11).Value = 6034.666499999999
$ws.Cells.Item(122, 12).Value = 8249.625
$ws.Cells.Item(122, 13).Value = -3584.666499999999
$ws.Cells.Item(122, 14).Value = -13149.625
# Row 132
$ws.Cells.Item(132, 8).Value = 45870.78
$ws.Cells.Item(132, 9).Value = 79404
$ws.Cells.Item(132, 10).Value = 2277.6
$ws.Cells.Item(132, 11).Value = 238212
$ws.Cells.Item(132, 12).Value = 6832.799999999999
$ws.Cells.Item(132, 13).Value = -235682
$ws.Cells.Item(132, 14).Value = -11892.8

$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 14).Value = $null
# Row 132
$ws.Cells.Item(132, 8).Value = 57310540
$ws.Cells.Item(132, 9).Value = 90401490
$ws.Cells.Item(132, 10).Value = 2158964.5
$ws.Cells.Item(132, 11).Value = 271204470
$ws.Cells.Item(132, 12).Value = 6476893.5
$ws.Cells.Item(132, 13).Value = -271201940
$ws.Cells.Item(132, 14).Value = -6481953.5

Write-Host "Applied all Sheet edits"